$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '67.258.85' },
    @{ Cell = 'E2'; Value = '  +0.82%  ' },
    @{ Cell = 'D3'; Value = '3.479.90' },
    @{ Cell = 'E3'; Value = '  -0.21%  ' },
    @{ Cell = 'D5'; Value = '593.53' },
    @{ Cell = 'E5'; Value = '  +0.21%  ' },
    @{ Cell = 'D6'; Value = '178.48' },
    @{ Cell = 'E6'; Value = '  +4.02%  ' },
    @{ Cell = 'E7'; Value = '  +0.01%  ' },
    @{ Cell = 'D8'; Value = '0.601' },
    @{ Cell = 'E8'; Value = '  +1.71%  ' },
    @{ Cell = 'D9'; Value = '3.481.24' },
    @{ Cell = 'E9'; Value = '  -0.15%  ' },
    @{ Cell = 'D10'; Value = '0.137' },
    @{ Cell = 'E10'; Value = '  +4.84%  ' },
    @{ Cell = 'D11'; Value = '7.07' },
    @{ Cell = 'E11'; Value = '  -2.34%  ' },
    @{ Cell = 'D12'; Value = '0.434' },
    @{ Cell = 'E12'; Value = '  +0.82%  ' },
    @{ Cell = 'D13'; Value = '4.080.34' },
    @{ Cell = 'E13'; Value = '  -0.27%  ' },
    @{ Cell = 'D14'; Value = '31.98' },
    @{ Cell = 'E14'; Value = '  +10.81%  ' },
    @{ Cell = 'E15'; Value = '  +1.69%  ' },
    @{ Cell = 'D16'; Value = '67.322.15' },
    @{ Cell = 'E16'; Value = '  +0.82%  ' },
    @{ Cell = 'E17'; Value = '  -0.26%  ' },
    @{ Cell = 'D18'; Value = '3.483.18' },
    @{ Cell = 'E18'; Value = '  -0.28%  ' },
    @{ Cell = 'E19'; Value = '  -0.08%  ' },
    @{ Cell = 'E20'; Value = '  +1.77%  ' },
    @{ Cell = 'D21'; Value = '388.76' },
    @{ Cell = 'E21'; Value = '  -0.89%  ' },
    @{ Cell = 'D22'; Value = '7.90' },
    @{ Cell = 'E22'; Value = '  -0.15%  ' },
    @{ Cell = 'D23'; Value = '73.97' },
    @{ Cell = 'E23'; Value = '  +1.68%  ' },
    @{ Cell = 'D24'; Value = '0.998' },
    @{ Cell = 'E24'; Value = '  -0.18%  ' },
    @{ Cell = 'D25'; Value = '0.536' },
    @{ Cell = 'E25'; Value = '  +0.29%  ' },
    @{ Cell = 'E26'; Value = '  +0.54%  ' },
    @{ Cell = 'E27'; Value = '  +0.78%  ' },
    @{ Cell = 'D28'; Value = '10.37' },
    @{ Cell = 'E28'; Value = '  +2.16%  ' },
    @{ Cell = 'E29'; Value = '  -3.14%  ' },
    @{ Cell = 'D30'; Value = '1.00' },
    @{ Cell = 'E30'; Value = '  +0.12%  ' },
    @{ Cell = 'D31'; Value = '6.16' },
    @{ Cell = 'E31'; Value = '  -0.26%  ' },
    @{ Cell = 'E32'; Value = '  -0.21%  ' },
    @{ Cell = 'E33'; Value = '  +0.52%  ' },
    @{ Cell = 'D34'; Value = '23.54' },
    @{ Cell = 'E34'; Value = '  -0.67%  ' },
    @{ Cell = 'D35'; Value = '7.35' },
    @{ Cell = 'E35'; Value = '  +0.34%  ' },
    @{ Cell = 'E37'; Value = '  -1.67%  ' },
    @{ Cell = 'D38'; Value = '164.10' },
    @{ Cell = 'E38'; Value = '  +0.61%  ' },
    @{ Cell = 'D39'; Value = '0.872' },
    @{ Cell = 'E39'; Value = '  -0.58%  ' },
    @{ Cell = 'E40'; Value = '  -0.47%  ' },
    @{ Cell = 'E41'; Value = '  +7.00%  ' },
    @{ Cell = 'E42'; Value = '  -0.13%  ' },
    @{ Cell = 'D43'; Value = '4.64' },
    @{ Cell = 'E43'; Value = '  -0.15%  ' },
    @{ Cell = 'D44'; Value = '2.844.26' },
    @{ Cell = 'E44'; Value = '  +1.65%  ' },
    @{ Cell = 'D45'; Value = '26.24' },
    @{ Cell = 'E45'; Value = '  +0.45%  ' },
    @{ Cell = 'D46'; Value = '26.94' },
    @{ Cell = 'E46'; Value = '  -0.60%  ' },
    @{ Cell = 'D47'; Value = '0.0720' },
    @{ Cell = 'E47'; Value = '  -2.53%  ' },
    @{ Cell = 'D48'; Value = '41.50' },
    @{ Cell = 'E48'; Value = '  -2.78%  ' },
    @{ Cell = 'D49'; Value = '0.0299' },
    @{ Cell = 'E49'; Value = '  -1.04%  ' },
    @{ Cell = 'D50'; Value = '336.08' },
    @{ Cell = 'E50'; Value = '  +0.18%  ' },
    @{ Cell = 'E51'; Value = '  -2.19%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
